{"js": "// Word JavaScript API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Reproduces the two meaningful, deterministic content changes captured by\n// the diff:\n//   1. A new, empty leading paragraph is inserted at the very start of the\n//      document body (before the first table), carrying paragraph-mark\n//      run formatting of rFonts/@w:cstheme=\"minorHAnsi\".\n//   2. The hidden \"_GoBack\" bookmark (Word's \"last edit location\" marker)\n//      is moved from the end of the document (after the \"El Camino de\n//      Santiago.\" bullet) to inside the \"...and math skills.\" sentence,\n//      right after \"math s\".\n\nconst body = context.document.body;\n\n// 1) Insert the new leading empty paragraph with explicit OOXML so the\n//    paragraph-mark run properties (rFonts cstheme=\"minorHAnsi\") match\n//    exactly, instead of the generic empty run insertParagraph() leaves.\nconst leadingParaOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:pPr><w:rPr><w:rFonts w:cstheme=\"minorHAnsi\"/></w:rPr></w:pPr></w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nbody.insertOoxml(leadingParaOoxml, \"Start\");\nawait context.sync();\n\n// 2) Relocate the \"_GoBack\" bookmark. Remove it from its old position\n//    (it is harmless if it is not there) and insert it right after\n//    \"math s\" inside the \"Strong computer science, graphics programming,\n//    and math skills.\" sentence.\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // No-op if the bookmark does not currently exist.\n}\n\nconst results = body.search(\"and math s\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const found = results.items[0];\n  const endOfMatch = found.getRange(\"End\");\n  endOfMatch.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is available as $d.\n#\n# Reproduces the two meaningful, deterministic content changes captured by\n# the diff:\n#   1. A new, empty leading paragraph is inserted at the very start of the\n#      document body (before the first table), carrying paragraph-mark\n#      run formatting of rFonts/@w:cstheme=\"minorHAnsi\".\n#   2. The hidden \"_GoBack\" bookmark (Word's \"last edit location\" marker)\n#      is moved from the end of the document (after the \"El Camino de\n#      Santiago.\" bullet) to inside the \"...and math skills.\" sentence,\n#      right after \"math s\".\n\n$d = $word.ActiveDocument\n\n# 1) Insert the new leading empty paragraph with explicit WordOpenXML so the\n#    paragraph-mark run properties (rFonts cstheme=\"minorHAnsi\") match\n#    exactly, instead of the generic empty run InsertParagraphBefore() leaves.\n$leadingParaXml = @\"\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:pPr><w:rPr><w:rFonts w:cstheme=\"minorHAnsi\"/></w:rPr></w:pPr></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$startRange = $d.Range(0, 0)\n$startRange.InsertXML($leadingParaXml)\n\n# 2) Relocate the \"_GoBack\" bookmark. Find \"and math s\" inside the\n#    \"Strong computer science, graphics programming, and math skills.\"\n#    sentence and (re)create the bookmark, collapsed, right after it.\n#    Adding a bookmark named \"_GoBack\" replaces any existing one of the\n#    same name, so the old location is automatically vacated.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"and math s\")\nif ($found) {\n    $bookmarkRange = $d.Range($findRange.End, $findRange.End)\n    $d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n}\n"}
